$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.989.74'
$ws.Range("E2").Value = '  -5.95%  '

# Row 3
$ws.Range("D3").Value = '3.344.16'
$ws.Range("E3").Value = '  -2.54%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("E5").Value = '  -3.24%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.97'
$ws.Range("E6").Value = '  +0.42%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").Value = '3.343.99'
$ws.Range("E8").Value = '  -2.55%  '

# Row 9
$ws.Range("E9").Value = '  -1.76%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.42'
$ws.Range("E10").Value = '  -2.38%  '

# Row 11
$ws.Range("E11").Value = '  -6.47%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.377'
$ws.Range("E12").Value = '  -1.70%  '

# Row 13
$ws.Range("D13").Value = '3.910.35'
$ws.Range("E13").Value = '  -2.60%  '

# Row 14
$ws.Range("E14").Value = '  -0.32%  '

# Row 15
$ws.Range("D15").Value = '3.332.84'
$ws.Range("E15").Value = '  -2.78%  '

# Row 16
$ws.Range("E16").Value = '  -5.35%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.70'
$ws.Range("E17").Value = '  -1.43%  '

# Row 18
$ws.Range("D18").Value = '60.172.29'
$ws.Range("E18").Value = '  -5.70%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.67'
$ws.Range("E19").Value = '  -0.17%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.43'
$ws.Range("E20").Value = '  +0.54%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.00'
$ws.Range("E21").Value = '  -8.84%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '354.23'
$ws.Range("E22").Value = '  -7.96%  '

# Row 23
$ws.Range("E23").Value = '  -1.28%  '

# Row 24
$ws.Range("D24").Value = '3.476.19'
$ws.Range("E24").Value = '  -2.58%  '

# Row 25
$ws.Range("E25").Value = '  -0.01%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.29'
$ws.Range("E26").Value = '  -6.34%  '

# Row 27
$ws.Range("E27").Value = '  +1.59%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.64'
$ws.Range("E28").Value = '  +15.15%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.49'
$ws.Range("E29").Value = '  +6.01%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.07%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.94'
$ws.Range("E31").Value = '  -0.38%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.154'
$ws.Range("E32").Value = '  -0.68%  '

# Row 33
$ws.Range("E33").Value = '  -4.31%  '

# Row 34
$ws.Range("E34").Value = '  -0.04%  '

# Row 35
$ws.Range("D35").Value = '3.372.94'
$ws.Range("E35").Value = '  -2.49%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.90'
$ws.Range("E36").Value = '  -0.16%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.38'
$ws.Range("E37").Value = '  +3.72%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.87'
$ws.Range("E38").Value = '  +1.15%  '

# Row 39
$ws.Range("E39").Value = '  -0.66%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '157.60'
$ws.Range("E40").Value = '  -3.93%  '

# Row 41
$ws.Range("E41").Value = '  -1.47%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.09%  '

# Row 43
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.38'
$ws.Range("E43").Value = '  +0.64%  '

# Row 44
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.20'
$ws.Range("E44").Value = '  +7.61%  '

# Row 45
$ws.Range("E45").Value = '  -5.03%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.69'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.72'
$ws.Range("E47").Value = '  +0.98%  '

# Row 48
$ws.Range("E48").Value = '  -3.24%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.83'
$ws.Range("E49").Value = '  +1.38%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.48'
$ws.Range("E50").Value = '  +10.31%  '

# Row 51
$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.41'
$ws.Range("E51").Value = '  +15.25%  '
